$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (names + totals), already sorted by total_registros descending,
# as found in the updated source file.
$names = @(
    "NEYRA PEREIRA MONICA AGNES ALEXANDRA",
    "MOGOLLON MORON KARIN LISSET",
    "TUANAMA PIZANGO ELIZABETH",
    "AGURTO TINEO CESIA JIMENA",
    "MORETO ESPINOZA CRISTIAN ESTEBAN",
    "SANCARRANCO SANCHEZ DE CRUZ GISSELA SHANI",
    "CHIROQUE YARLEQUE BETTY ELIZABETH",
    "HERRERA JUAN MANUEL",
    "SILVA ALVARADO EVELYN DE JESUS",
    "TENE TRABUCCO GIAN PIERRE"
)

$totals = @(139, 138, 113, 113, 107, 104, 102, 102, 96, 92)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $totals[$i]
}
